$d = $word.ActiveDocument

# Helper: replace the word "Columns" with "Rows" inside a run that is found via
# a unique surrounding context string (e.g. "ds.Columns.TheDouble"), while
# preserving the surrounding runs (so they do not get merged together).
# The trick: toggle Bold on the "Columns" sub-range first so it becomes a
# separate run (different formatting than its neighbours), then replace the
# text (this does not touch the neighbouring runs any more), then toggle
# Bold back off on the (now shorter) replaced range so the formatting is
# restored to the original (identical) formatting - the run boundaries
# created earlier remain intact because only direct text edits cause the
# surrounding runs to be coalesced.

function Replace-ColumnsToken([string]$context) {
    $rng = $d.Content
    $found = $rng.Find.Execute($context, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "context not found:" $context
        return
    }
    $ctxStart = $rng.Start
    $ctxEnd = $rng.End

    $inner = $d.Range($ctxStart, $ctxEnd)
    $innerFound = $inner.Find.Execute("Columns", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $innerFound) {
        Write-Host "Columns not found within context:" $context
        return
    }
    $tokenStart = $inner.Start
    $tokenEnd = $inner.End

    $tokenRange = $d.Range($tokenStart, $tokenEnd)
    $tokenRange.Font.Bold = 1

    $tokenRange2 = $d.Range($tokenStart, $tokenEnd)
    $tokenRange2.Text = "Rows"

    $newEnd = $tokenStart + 4
    $tokenRange3 = $d.Range($tokenStart, $newEnd)
    $tokenRange3.Font.Bold = 0
}

Replace-ColumnsToken "ds.Columns.TheDouble"
Replace-ColumnsToken "ds.Columns.TheDate"
